# Inserts a new weekly price record at row 203 of the Coliflor sheet,
# pushing the existing rows 203-335 down to 204-336 (dimension grows
# from A1:R335 to A1:R336).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 203..335 down by one row.
$ws.Rows.Item(203).Insert()

# Populate the newly inserted row 203 with the new record.
$ws.Cells.Item(203, 1).Value  = 10
$ws.Cells.Item(203, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(203, 3).Value  = "La Araucanía"
$ws.Cells.Item(203, 4).Value  = 44606
$ws.Cells.Item(203, 5).Value  = 9
$ws.Cells.Item(203, 6).Value  = 100112008
$ws.Cells.Item(203, 7).Value  = "Coliflor"
$ws.Cells.Item(203, 8).Value  = "Sin especificar"
$ws.Cells.Item(203, 9).Value  = "Primera"
$ws.Cells.Item(203, 10).Value = 200
$ws.Cells.Item(203, 11).Value = 1200
$ws.Cells.Item(203, 12).Value = 1200
$ws.Cells.Item(203, 13).Value = 1200
$ws.Cells.Item(203, 14).Value = "`$/unidad"
$ws.Cells.Item(203, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(203, 16).Value = 1200
$ws.Cells.Item(203, 17).Value = 1
$ws.Cells.Item(203, 18).Value = "Hortaliza"
